$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 60.8
$ws.Range("K3").Value = 55
$ws.Range("K4").Value = 51.2
$ws.Range("K5").Value = 48.8

$ws.Range("N2").Value = 85.96878041621773
$ws.Range("N3").Value = 85.96878041621773
$ws.Range("N4").Value = 85.96878041621773
$ws.Range("N5").Value = 85.96878041621773
